$d = $word.ActiveDocument
$c = $d.Content

# Rework the "Du hast Krallen gefunden! ..." tutorial text: move the "\n"
# line-break markers, add a couple of commas, and join "Inventar zu
# öffnen"/"und wähle" without the old separator.
$c.Find.Execute("\n Jetzt kann das Eichhörnchen\n", $false, $false, $false, $false, $false, $true, 1, $false, "Jetzt kann das \n Eichhörnchen", 2)
$c.Find.Execute(" \n um dich daran festzuhalten. \n", $false, $false, $false, $false, $false, $true, 1, $false, ", um dich daran \n festzuhalten.", 2)
$c.Find.Execute("bzw. Shift rechts \n gedrückt", $false, $false, $false, $false, $false, $true, 1, $false, "bzw.  Shift rechts gedrückt,", 2)
$c.Find.Execute(" \n und wähle", $false, $false, $false, $false, $false, $true, 1, $false, "und wähle", 2)
$c.Find.Execute("bzw. WASD ein Item aus.", $false, $false, $false, $false, $false, $true, 1, $false, "bzw. \n WASD ein Item aus.", 2)

# Move the "_GoBack" bookmark from the "2. Level" chest text onto this
# paragraph (spanning everything except the leading curly quote).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("„Du hast Krallen")) {
        $rng = $d.Range($p.Range.Start + 1, $p.Range.End)
        $d.Bookmarks.Add("_GoBack", $rng)
        break
    }
}
